# edit.ps1 - apply "slight edits on parsers.pptx" changes via PowerPoint COM-interop
$p = $ppt.ActivePresentation

# --- Slide 11 ("Implementation"): Content Placeholder 2, paragraph 3 ---
# "Parsing K itself becomes special case: a two-level parser (OUTER, then KAST)"
# -> "Parsing K itself becomes special case: a two-level parser (outer, then inner)"
# split across five runs
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(2)
$para11 = $sh11.TextFrame.TextRange.Paragraphs(3)
$para11.Text = "Parsing K itself becomes special case: a two-level parser "
[void]$para11.InsertAfter("(outer, ")
[void]$para11.InsertAfter("then ")
[void]$para11.InsertAfter("inner")
[void]$para11.InsertAfter(")")

# --- Slide 13: Content Placeholder 8, paragraph 1 ---
# "(2) Define BUBBLE with sort <b>Bubble</b> as a list of <b>BubbleItem</b> elements:"
# -> "(2) Define BUBBLE with sort  <b>Bubble</b>   as a list of  <b>BubbleItem</b>  elements:"
$s13 = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(4)
$para13 = $sh13.TextFrame.TextRange.Paragraphs(1)
$para13.Text = "(2) Define BUBBLE with "
[void]$para13.InsertAfter("sort  ")
[void]$para13.InsertAfter("Bubble")
[void]$para13.InsertAfter(" ")
[void]$para13.InsertAfter(" as ")
[void]$para13.InsertAfter("a list ")
[void]$para13.InsertAfter("of  ")
[void]$para13.InsertAfter("BubbleItem")
[void]$para13.InsertAfter(" ")
[void]$para13.InsertAfter(" elements")
[void]$para13.InsertAfter(":")
# re-apply bold/Courier New formatting to the two code tokens (lost their
# run-level formatting when the paragraph text was rebuilt above)
$bubble1 = $para13.Characters(30, 6)
$bubble1.Font.Bold = $true
$bubble1.Font.Name = "Courier New"
$bubble2 = $para13.Characters(52, 10)
$bubble2.Font.Bold = $true
$bubble2.Font.Name = "Courier New"

# --- Slide 2: three textboxes, merge the multiple runs in each paragraph into one ---
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(8).TextFrame.TextRange.Paragraphs(1).Text = "=   K syntax"
$s2.Shapes.Item(10).TextFrame.TextRange.Paragraphs(1).Text = "=   PL/library concrete syntax"
$s2.Shapes.Item(12).TextFrame.TextRange.Paragraphs(1).Text = "=   KAST (PL abstract syntax)"

# --- Slide 4: Content Placeholder 2, paragraph 7 ---
# "Use fast KAST parser when bubble contains PL syntax only using KAST"
# -> "Use fast KAST parser when bubble contains PL syntax using KAST only"
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$para4 = $sh4.TextFrame.TextRange.Paragraphs(7)
$para4.Text = "Use fast KAST parser when bubble contains PL syntax "
[void]$para4.InsertAfter("using KAST ")
[void]$para4.InsertAfter("only")
